$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price (column D) and 1h volume-change (column E)
# figures for this data refresh. Some of the new "Price" strings look like
# plain decimal numbers (e.g. "228.87"), so a leading apostrophe is used to
# force Excel to store them as text (matching the sheet's existing
# text-stored price strings such as "38.136.64") instead of silently
# auto-converting them to numeric values.
$ws.Range("D2").Value = '38.136.64'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").Value = '2.092.36'
$ws.Range("E3").Value = '  +2.72%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''228.87'
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").Value = '''0.612'
$ws.Range("E6").Value = '  +0.29%  '

$ws.Range("D7").Value = '''60.66'
$ws.Range("E7").Value = '  -0.55%  '

$ws.Range("D9").Value = '''0.379'
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("E10").Value = '  +3.06%  '

$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").Value = '2.400.66'
$ws.Range("E12").Value = '  +2.64%  '

$ws.Range("D13").Value = '''14.63'
$ws.Range("E13").Value = '  +0.70%  '

$ws.Range("D14").Value = '''22.15'
$ws.Range("E14").Value = '  +3.13%  '

$ws.Range("E15").Value = '  +5.62%  '

$ws.Range("D16").Value = '''0.774'
$ws.Range("E16").Value = '  +1.52%  '

$ws.Range("D17").Value = '2.111.11'
$ws.Range("E17").Value = '  +4.05%  '

$ws.Range("D18").Value = '38.060.58'
$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("D19").Value = '''6.01'
$ws.Range("E19").Value = '  +1.70%  '

$ws.Range("D20").Value = '''70.11'
$ws.Range("E20").Value = '  +0.42%  '

$ws.Range("D21").Value = '0.0₃0835'
$ws.Range("E21").Value = '  +1.26%  '

$ws.Range("D22").Value = '''223.95'
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("E24").Value = '  +1.26%  '

$ws.Range("E25").Value = '  +3.00%  '

$ws.Range("D26").Value = '''169.73'
$ws.Range("E26").Value = '  +1.33%  '

$ws.Range("D27").Value = '''9.42'
$ws.Range("E27").Value = '  +0.79%  '

$ws.Range("E28").Value = '  +0.07%  '

$ws.Range("D29").Value = '''18.95'
$ws.Range("E29").Value = '  +0.26%  '

$ws.Range("D30").Value = '''1.35'
$ws.Range("E30").Value = '  +5.80%  '

$ws.Range("E31").Value = '  -0.48%  '

$ws.Range("E32").Value = '  +5.30%  '

$ws.Range("E33").Value = '  +4.04%  '

$ws.Range("D34").Value = '''4.43'
$ws.Range("E34").Value = '  +0.23%  '

$ws.Range("D35").Value = '''0.0606'
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").Value = '''2.40'
$ws.Range("E36").Value = '  +4.33%  '

$ws.Range("D37").Value = '''6.40'
$ws.Range("E37").Value = '  +0.34%  '

$ws.Range("D38").Value = '''3.52'
$ws.Range("E38").Value = '  +5.01%  '

$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("E40").Value = '  +1.93%  '

$ws.Range("D41").Value = '1.557.13'
$ws.Range("E41").Value = '  +1.02%  '

$ws.Range("D42").Value = '''100.17'
$ws.Range("E42").Value = '  +3.89%  '

$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("E44").Value = '  +0.91%  '

$ws.Range("D45").Value = '''0.0912'
$ws.Range("E45").Value = '  -0.29%  '

$ws.Range("D46").Value = '''4.14'
$ws.Range("E46").Value = '  +3.16%  '

$ws.Range("E47").Value = '  +0.88%  '

$ws.Range("D48").Value = '''7.49'
$ws.Range("E48").Value = '  +5.65%  '

$ws.Range("E49").Value = '  +1.23%  '

$ws.Range("D50").Value = '''2.99'
$ws.Range("E50").Value = '  +0.78%  '

$ws.Range("D51").Value = '2.288.12'
$ws.Range("E51").Value = '  +2.70%  '
